$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.637.13'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '2.472.95'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.552'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.84%  '
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0868'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '33.17'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.51%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '2.854.55'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.56'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.87%  '
$ws.Range('D16').Value = '2.474.00'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('E17').Value = '  +3.09%  '
$ws.Range('D18').Value = '41.586.27'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('D20').Value = '0.0₃0945'
$ws.Range('E20').Value = '  +1.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '70.71'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.21'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.51%  '
$ws.Range('E24').Value = '  +2.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.95'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.54%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.90%  '
$ws.Range('E28').Value = '  -0.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.70'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.55'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '156.68'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.92%  '
$ws.Range('E32').Value = '  +0.84%  '
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('E34').Value = '  +1.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.56'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.30%  '
$ws.Range('B37').Value = 'Stellar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.116'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.90%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.84'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.44%  '
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('E40').Value = '  +2.61%  '
$ws.Range('E41').Value = '  +2.39%  '
$ws.Range('E42').Value = '  +3.89%  '
$ws.Range('D43').Value = '1.986.56'
$ws.Range('E43').Value = '  +1.00%  '
$ws.Range('E44').Value = '  +0.91%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.80'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.66%  '
$ws.Range('E46').Value = '  +3.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.86%  '
$ws.Range('D48').Value = '2.714.37'
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.89'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.64'
$ws.Range('D50').Style = 'Normal'
